$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Bump the "pt" (receipt) status of the existing rows from 2 to 3 ---
$ws.Range("C2:C13").Value = 3

# --- 2) Copy the cell formatting used by the existing data rows for column B
#        (date, style index 1) down onto the new rows so the new date cells
#        carry the same "m/d/yyyy" number format as the rest of the table,
#        without introducing any new style records. ---
$ws.Range("B4").Copy()
$ws.Range("B14:B25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3) Append a duplicate batch of rows (14-25) for the next day
#        (11/14/2020, serial 44056) under receipt #3, reusing the same
#        po/bmi_code/qty values as the original rows 2-13. ---
$poValues   = @(6800083955,6800083955,6800083955,6800083954,6800083955,6800083955,6800083955,6800083955,6800083955,6800083955,6800083955,6800083955)
$codeValues = @("TA23D7GANV6","TA40Y3LANC2","TA5762NANV1","TA10VG3ANV1","TA18410ANO1","TA23D7GANV6","TA18410ANO1","TA11J0XANY2","TA10UF7ANH2","TA18410ANO1","TA194G3ANV1","TA23D7GANV6")
$qtyValues  = @(2198.268,394.98,140.74,27.24,608.36,155.268,1198.56,4.54,690.07999999999993,349.58,631.05999999999995,650.74400000000003)

for ($i = 0; $i -lt 12; $i++) {
    $r = 14 + $i
    $ws.Range("A$r").Value = $poValues[$i]
    $ws.Range("B$r").Value = 44056
    $ws.Range("C$r").Value = 3
    $ws.Range("D$r").Value = $codeValues[$i]
    $ws.Range("E$r").Value = $qtyValues[$i]
}

# --- 4) Grow Table1 to cover the newly added rows. ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E25"))

# --- 5) Match the author's final selection. ---
$ws.Range("E9").Select()
